# Informe sumo.docx - "En cuanto al montaje de los sensores infrarrojos
# Sharp ..." paragraph: replace "dentro del" with "para el", and add a new
# paragraph about the state-machine programming right after it (before the
# _GoBack bookmark, which stays anchored to the end of the original
# paragraph).

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive text instead of a hard
# coded index, so the script is resilient to any paragraph numbering
# differences.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "En cuanto al montaje de los sensores infrarrojos Sharp*") {
        $target = $para
        break
    }
}

$paraText = $target.Range.Text

# --- Step 1: turn "... funcionamiento dentro del sistema." into
#             "... funcionamiento para el sistema." -----------------------
$oldPhrase = "dentro del "
$newPhrase = "para el "
$idx = $paraText.IndexOf($oldPhrase)
$absStart = $target.Range.Start + $idx

$rOld = $d.Range($absStart, $absStart + $oldPhrase.Length)
$rOld.Text = $newPhrase

# Re-assert (no-op) character formatting on the freshly typed text so it is
# kept as its own run instead of being silently coalesced with its
# neighbours when the package is serialized - mirrors what happens in Word
# when text is typed over a selection.
$rNew = $d.Range($absStart, $absStart + $newPhrase.Length)
$rNew.Bold = 1
$rNew.Bold = 0

# Likewise re-assert formatting over the untouched remainder of the
# paragraph (from "sensores ubicados ..." onward) so that its existing run
# boundary survives the save instead of merging into the first run.
$paraText2 = $target.Range.Text
$idxB = $paraText2.IndexOf("sensores ubicados en la parte inferior")
if ($idxB -ge 0) {
    $absB = $target.Range.Start + $idxB
    $endOfPara = $target.Range.End - 1
    $rB = $d.Range($absB, $endOfPara)
    $rB.Bold = 1
    $rB.Bold = 0
}

# --- Step 2: insert the new paragraph right after the edited one, before
#             the trailing _GoBack bookmark. -------------------------------
$endOfPara2 = $target.Range.End
$insertionPoint = $d.Range($endOfPara2 - 1, $endOfPara2 - 1)
$insertionPoint.InsertParagraphAfter()

$newParaIndex = $target.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newParaText = "Para la programación, siendo esta en base a máquinas de estado primero se cargó un programa piloto la el control de los motores, de esta manera verificando su funcionamiento, dicho programa piloto mediante máquinas de est5ado, genera lo necesario para poder ir hacia adelante, izquierda, derecha y atrás, así mismo el estado de avanzar y de frenado, para el control básico de los motores así verificando que el circuito y todos sus componentes funcionan de manera apropiada"
$newPara.Range.InsertAfter($newParaText)
